$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2129.8096
$ws.Range("I40").Value = 909.1429000000001
$ws.Range("K40").Value = 909.1429000000001
$ws.Range("M40").Value = -734.1429000000001
$ws.Range("H74").Value = 8940.076999999999
$ws.Range("I74").Value = 5152.625
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 5152.625
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -4216.625
$ws.Range("N74").Value = -16872
$ws.Range("H76").Value = 4256.6665
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 4885
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 4885
$ws.Range("M76").Value = -2685
$ws.Range("N76").Value = -5515
$ws.Range("H77").Value = 8940.076999999999
$ws.Range("I77").Value = 5152.625
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 25763.125
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -21083.125
$ws.Range("N77").Value = -84360
$ws.Range("H79").Value = 4256.6665
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 4885
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 4885
$ws.Range("M79").Value = -1908
$ws.Range("N79").Value = -7069
$ws.Range("H93").Value = 60100
$ws.Range("J93").Value = 60100
$ws.Range("L93").Value = 60100
$ws.Range("N93").Value = -65092
$ws.Range("H100").Value = 5345.364
$ws.Range("I100").Value = 3133.1667
$ws.Range("K100").Value = 3133.1667
$ws.Range("M100").Value = -2592.1667
$ws.Range("H106").Value = 3849.6667
$ws.Range("I106").Value = 3849.6667
$ws.Range("K106").Value = 3849.6667
$ws.Range("M106").Value = -3218.6667
$ws.Range("H116").Value = 14948
$ws.Range("I116").Value = 14932.333
$ws.Range("K116").Value = 14932.333
$ws.Range("M116").Value = -11490.333
$ws.Range("H137").Value = 6839.2383
$ws.Range("I137").Value = 8317.733
$ws.Range("J137").Value = 3143
$ws.Range("K137").Value = 24953.199
$ws.Range("L137").Value = 9429
$ws.Range("M137").Value = -22403.199
$ws.Range("N137").Value = -14529
$ws.Range("H141").Value = 3897.5
$ws.Range("I141").Value = 3795
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 11385
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -6205
$ws.Range("N141").Value = -22360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2251.5
$ws.Range("I45").Value = 1389.875
$ws.Range("J45").Value = 3974.75
$ws.Range("K45").Value = 1389.875
$ws.Range("L45").Value = 3974.75
$ws.Range("M45").Value = -1012.875
$ws.Range("N45").Value = -4728.75
$ws.Range("H62").Value = 99998
$ws.Range("J62").Value = 99998
$ws.Range("L62").Value = 99998
$ws.Range("N62").Value = -101246
$ws.Range("H63").Value = 7604.778
$ws.Range("I63").Value = 7183.625
$ws.Range("J63").Value = 7941.7
$ws.Range("K63").Value = 7183.625
$ws.Range("L63").Value = 7941.7
$ws.Range("M63").Value = -6497.625
$ws.Range("N63").Value = -9313.700000000001
$ws.Range("H65").Value = 99998
$ws.Range("J65").Value = 99998
$ws.Range("L65").Value = 299994
$ws.Range("N65").Value = -306234
$ws.Range("H66").Value = 7604.778
$ws.Range("I66").Value = 7183.625
$ws.Range("J66").Value = 7941.7
$ws.Range("K66").Value = 35918.125
$ws.Range("L66").Value = 39708.5
$ws.Range("M66").Value = -32486.125
$ws.Range("N66").Value = -46572.5
$ws.Range("H74").Value = 5421.5454
$ws.Range("I74").Value = 4752.385
$ws.Range("K74").Value = 4752.385
$ws.Range("M74").Value = -3878.385
$ws.Range("H77").Value = 5421.5454
$ws.Range("I77").Value = 4752.385
$ws.Range("K77").Value = 23761.925
$ws.Range("M77").Value = -19393.925

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 388
$ws.Range("J64").Value = 388
$ws.Range("L64").Value = 388
$ws.Range("N64").Value = -838
$ws.Range("H67").Value = 388
$ws.Range("J67").Value = 388
$ws.Range("L67").Value = 388
$ws.Range("N67").Value = -1948
$ws.Range("H86").Value = 50002910
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 50002910
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = 50002910
$ws.Range("N86").Value = -50005156
$ws.Range("H89").Value = 50002910
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 50002910
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = 250014550
$ws.Range("N89").Value = -250025782
$ws.Range("H130").Value = 289997
$ws.Range("J130").Value = 289997
$ws.Range("L130").Value = 289997
$ws.Range("N130").Value = -300037
$ws.Range("H131").Value = 232196.89
$ws.Range("J131").Value = 227471.62
$ws.Range("L131").Value = 227471.62
$ws.Range("N131").Value = -237551.62

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 12583.333
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 14900
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 14900
$ws.Range("M23").Value = -760
$ws.Range("N23").Value = -15380
$ws.Range("H27").Value = 12583.333
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 14900
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 14900
$ws.Range("M27").Value = -808
$ws.Range("N27").Value = -15284
$ws.Range("H31").Value = 2441.5588
$ws.Range("J31").Value = 4201.357
$ws.Range("L31").Value = 4201.357
$ws.Range("N31").Value = -4791.357
$ws.Range("H34").Value = 2441.5588
$ws.Range("J34").Value = 4201.357
$ws.Range("L34").Value = 4201.357
$ws.Range("N34").Value = -4605.357
$ws.Range("H62").Value = 6914
$ws.Range("I62").Value = 5502.5
$ws.Range("J62").Value = 7478.6
$ws.Range("K62").Value = 5502.5
$ws.Range("L62").Value = 7478.6
$ws.Range("M62").Value = -4878.5
$ws.Range("N62").Value = -8726.6
$ws.Range("H65").Value = 6914
$ws.Range("I65").Value = 5502.5
$ws.Range("J65").Value = 7478.6
$ws.Range("K65").Value = 27512.5
$ws.Range("L65").Value = 37393
$ws.Range("M65").Value = -24392.5
$ws.Range("N65").Value = -43633

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2917.2727
$ws.Range("I107").Value = 2441.7144
$ws.Range("K107").Value = 7325.1432
$ws.Range("M107").Value = -5405.1432
$ws.Range("H111").Value = 22375
$ws.Range("I111").Value = 812.5
$ws.Range("K111").Value = 2437.5
$ws.Range("M111").Value = 629.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2073.3845
$ws.Range("I126").Value = 2054.5
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 6163.5
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -3693.5
$ws.Range("N126").Value = -11840

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 61362.79
$ws.Range("I16").Value = 81971.07000000001
$ws.Range("J16").Value = 3659.6
$ws.Range("K16").Value = 81971.07000000001
$ws.Range("L16").Value = 3659.6
$ws.Range("M16").Value = -81801.07000000001
$ws.Range("N16").Value = -3999.6
$ws.Range("H40").Value = 7605.926
$ws.Range("I40").Value = 6885.8823
$ws.Range("K40").Value = 6885.8823
$ws.Range("M40").Value = -6749.8823
$ws.Range("H61").Value = 5633.579
$ws.Range("I61").Value = 2219
$ws.Range("J61").Value = 7625.4165
$ws.Range("K61").Value = 2219
$ws.Range("L61").Value = 7625.4165
$ws.Range("M61").Value = -2017
$ws.Range("N61").Value = -8029.4165
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H68").Value = 5893.3335
$ws.Range("I68").Value = 3155.5557
$ws.Range("K68").Value = 3155.5557
$ws.Range("M68").Value = -2406.5557
$ws.Range("H71").Value = 5893.3335
$ws.Range("I71").Value = 3155.5557
$ws.Range("K71").Value = 15777.7785
$ws.Range("M71").Value = -12033.7785
$ws.Range("H113").Value = 5633.579
$ws.Range("I113").Value = 2219
$ws.Range("J113").Value = 7625.4165
$ws.Range("K113").Value = 2219
$ws.Range("L113").Value = 7625.4165
$ws.Range("M113").Value = -49
$ws.Range("N113").Value = -11965.4165
$ws.Range("H132").Value = 5297.9644
$ws.Range("I132").Value = 5016.3335
$ws.Range("J132").Value = 6142.857
$ws.Range("K132").Value = 15049.0005
$ws.Range("L132").Value = 18428.571
$ws.Range("M132").Value = -12519.0005
$ws.Range("N132").Value = -23488.571
$ws.Range("H136").Value = 3514.2144
$ws.Range("I136").Value = 3169.3845
$ws.Range("J136").Value = 7997
$ws.Range("K136").Value = 9508.1535
$ws.Range("L136").Value = 23991
$ws.Range("M136").Value = -6958.1535
$ws.Range("N136").Value = -29091

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10500
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 11333.333
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 11333.333
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -12581.333
$ws.Range("H65").Value = 10500
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 11333.333
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 56666.665
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -62906.665
$ws.Range("H126").Value = 3316.182
$ws.Range("I126").Value = 3247.8
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 9743.400000000001
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -7273.400000000001
$ws.Range("N126").Value = -16940
$ws.Range("H136").Value = 2532.4
$ws.Range("I136").Value = 2209.9333
$ws.Range("K136").Value = 6629.7999
$ws.Range("M136").Value = -4079.7999
